$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = 44838
$ws.Cells.Item(2, 9).Value = 'Primera'
$ws.Cells.Item(2, 10).Value = 120
$ws.Cells.Item(2, 11).Value = 6500
$ws.Cells.Item(2, 12).Value = 7000
$ws.Cells.Item(2, 13).Value = 6750
$ws.Cells.Item(2, 15).Value = 'Provincia de Diguillín'
$ws.Cells.Item(2, 16).Value = 675

# Row 3
$ws.Cells.Item(3, 4).Value = 44810
$ws.Cells.Item(3, 9).Value = 'Primera'
$ws.Cells.Item(3, 10).Value = 60
$ws.Cells.Item(3, 11).Value = 7000
$ws.Cells.Item(3, 12).Value = 8000
$ws.Cells.Item(3, 13).Value = 7500
$ws.Cells.Item(3, 15).Value = 'Provincia de Diguillín'
$ws.Cells.Item(3, 16).Value = 750

# Row 4
$ws.Cells.Item(4, 4).Value = 44846
$ws.Cells.Item(4, 9).Value = 'Primera'
$ws.Cells.Item(4, 10).Value = 100
$ws.Cells.Item(4, 11).Value = 6500
$ws.Cells.Item(4, 12).Value = 7000
$ws.Cells.Item(4, 13).Value = 6750
$ws.Cells.Item(4, 15).Value = 'Provincia de Diguillín'
$ws.Cells.Item(4, 16).Value = 675

# Row 5
$ws.Cells.Item(5, 4).Value = 44841
$ws.Cells.Item(5, 9).Value = 'Primera'
$ws.Cells.Item(5, 10).Value = 60
$ws.Cells.Item(5, 11).Value = 6500
$ws.Cells.Item(5, 12).Value = 7000
$ws.Cells.Item(5, 13).Value = 6750
$ws.Cells.Item(5, 15).Value = 'Provincia de Diguillín'
$ws.Cells.Item(5, 16).Value = 675

# Row 6
$ws.Cells.Item(6, 4).Value = 44798
$ws.Cells.Item(6, 9).Value = 'Primera'
$ws.Cells.Item(6, 10).Value = 80
$ws.Cells.Item(6, 11).Value = 7000
$ws.Cells.Item(6, 12).Value = 7000
$ws.Cells.Item(6, 13).Value = 7000
$ws.Cells.Item(6, 15).Value = 'Provincia de Diguillín'
$ws.Cells.Item(6, 16).Value = 700

# Row 7
$ws.Cells.Item(7, 4).Value = 44784
$ws.Cells.Item(7, 9).Value = 'Primera'
$ws.Cells.Item(7, 10).Value = 100
$ws.Cells.Item(7, 11).Value = 8000
$ws.Cells.Item(7, 12).Value = 9000
$ws.Cells.Item(7, 13).Value = 8500
$ws.Cells.Item(7, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(7, 16).Value = 850

# Row 8
$ws.Cells.Item(8, 4).Value = 44817
$ws.Cells.Item(8, 9).Value = 'Primera'
$ws.Cells.Item(8, 10).Value = 60
$ws.Cells.Item(8, 11).Value = 7000
$ws.Cells.Item(8, 12).Value = 7000
$ws.Cells.Item(8, 13).Value = 7000
$ws.Cells.Item(8, 15).Value = 'Provincia de Diguillín'
$ws.Cells.Item(8, 16).Value = 700

# Row 9
$ws.Cells.Item(9, 4).Value = 44817
$ws.Cells.Item(9, 9).Value = 'Segunda'
$ws.Cells.Item(9, 10).Value = 60
$ws.Cells.Item(9, 11).Value = 8000
$ws.Cells.Item(9, 12).Value = 8000
$ws.Cells.Item(9, 13).Value = 8000
$ws.Cells.Item(9, 15).Value = 'Provincia de Diguillín'
$ws.Cells.Item(9, 16).Value = 800

# Row 10
$ws.Cells.Item(10, 4).Value = 44804
$ws.Cells.Item(10, 9).Value = 'Primera'
$ws.Cells.Item(10, 10).Value = 80
$ws.Cells.Item(10, 11).Value = 7000
$ws.Cells.Item(10, 12).Value = 7500
$ws.Cells.Item(10, 13).Value = 7250
$ws.Cells.Item(10, 15).Value = 'Provincia de Diguillín'
$ws.Cells.Item(10, 16).Value = 725

# Row 11
$ws.Cells.Item(11, 4).Value = 44790
$ws.Cells.Item(11, 9).Value = 'Primera'
$ws.Cells.Item(11, 10).Value = 60
$ws.Cells.Item(11, 11).Value = 8500
$ws.Cells.Item(11, 12).Value = 9000
$ws.Cells.Item(11, 13).Value = 8750
$ws.Cells.Item(11, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(11, 16).Value = 875

# Row 12
$ws.Cells.Item(12, 4).Value = 44203
$ws.Cells.Item(12, 9).Value = 'Primera'
$ws.Cells.Item(12, 10).Value = 27
$ws.Cells.Item(12, 11).Value = 7000
$ws.Cells.Item(12, 12).Value = 8000
$ws.Cells.Item(12, 13).Value = 7556
$ws.Cells.Item(12, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(12, 16).Value = 756

# Row 13
$ws.Cells.Item(13, 4).Value = 44812
$ws.Cells.Item(13, 9).Value = 'Primera'
$ws.Cells.Item(13, 10).Value = 60
$ws.Cells.Item(13, 11).Value = 7000
$ws.Cells.Item(13, 12).Value = 8000
$ws.Cells.Item(13, 13).Value = 7500
$ws.Cells.Item(13, 15).Value = 'Provincia de Diguillín'
$ws.Cells.Item(13, 16).Value = 750

# Row 14
$ws.Cells.Item(14, 4).Value = 44775
$ws.Cells.Item(14, 9).Value = 'Primera'
$ws.Cells.Item(14, 10).Value = 60
$ws.Cells.Item(14, 11).Value = 8000
$ws.Cells.Item(14, 12).Value = 8000
$ws.Cells.Item(14, 13).Value = 8000
$ws.Cells.Item(14, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(14, 16).Value = 800

# Row 15
$ws.Cells.Item(15, 4).Value = 44782
$ws.Cells.Item(15, 9).Value = 'Primera'
$ws.Cells.Item(15, 10).Value = 120
$ws.Cells.Item(15, 11).Value = 8000
$ws.Cells.Item(15, 12).Value = 9000
$ws.Cells.Item(15, 13).Value = 8500
$ws.Cells.Item(15, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(15, 16).Value = 850

# Row 16
$ws.Cells.Item(16, 4).Value = 44799
$ws.Cells.Item(16, 9).Value = 'Primera'
$ws.Cells.Item(16, 10).Value = 60
$ws.Cells.Item(16, 11).Value = 7000
$ws.Cells.Item(16, 12).Value = 7000
$ws.Cells.Item(16, 13).Value = 7000
$ws.Cells.Item(16, 15).Value = 'Provincia de Diguillín'
$ws.Cells.Item(16, 16).Value = 700

# Row 17
$ws.Cells.Item(17, 4).Value = 44791
$ws.Cells.Item(17, 9).Value = 'Primera'
$ws.Cells.Item(17, 10).Value = 100
$ws.Cells.Item(17, 11).Value = 8500
$ws.Cells.Item(17, 12).Value = 9000
$ws.Cells.Item(17, 13).Value = 8750
$ws.Cells.Item(17, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(17, 16).Value = 875

# Row 18
$ws.Cells.Item(18, 4).Value = 44831
$ws.Cells.Item(18, 9).Value = 'Primera'
$ws.Cells.Item(18, 10).Value = 60
$ws.Cells.Item(18, 11).Value = 7000
$ws.Cells.Item(18, 12).Value = 7500
$ws.Cells.Item(18, 13).Value = 7250
$ws.Cells.Item(18, 15).Value = 'Provincia de Diguillín'
$ws.Cells.Item(18, 16).Value = 725

# Row 19
$ws.Cells.Item(19, 4).Value = 44211
$ws.Cells.Item(19, 9).Value = 'Primera'
$ws.Cells.Item(19, 10).Value = 28
$ws.Cells.Item(19, 11).Value = 8000
$ws.Cells.Item(19, 12).Value = 8500
$ws.Cells.Item(19, 13).Value = 8214
$ws.Cells.Item(19, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(19, 16).Value = 821

# Row 20
$ws.Cells.Item(20, 4).Value = 44847
$ws.Cells.Item(20, 9).Value = 'Primera'
$ws.Cells.Item(20, 10).Value = 100
$ws.Cells.Item(20, 11).Value = 6500
$ws.Cells.Item(20, 12).Value = 7000
$ws.Cells.Item(20, 13).Value = 6750
$ws.Cells.Item(20, 15).Value = 'Provincia de Diguillín'
$ws.Cells.Item(20, 16).Value = 675

# Row 21
$ws.Cells.Item(21, 4).Value = 44813
$ws.Cells.Item(21, 9).Value = 'Primera'
$ws.Cells.Item(21, 10).Value = 120
$ws.Cells.Item(21, 11).Value = 7000
$ws.Cells.Item(21, 12).Value = 7500
$ws.Cells.Item(21, 13).Value = 7250
$ws.Cells.Item(21, 15).Value = 'Provincia de Diguillín'
$ws.Cells.Item(21, 16).Value = 725

# Row 22
$ws.Cells.Item(22, 4).Value = 44819
$ws.Cells.Item(22, 9).Value = 'Primera'
$ws.Cells.Item(22, 10).Value = 100
$ws.Cells.Item(22, 11).Value = 7000
$ws.Cells.Item(22, 12).Value = 8000
$ws.Cells.Item(22, 13).Value = 7500
$ws.Cells.Item(22, 15).Value = 'Provincia de Diguillín'
$ws.Cells.Item(22, 16).Value = 750

# Row 23
$ws.Cells.Item(23, 4).Value = 44806
$ws.Cells.Item(23, 9).Value = 'Primera'
$ws.Cells.Item(23, 10).Value = 120
$ws.Cells.Item(23, 11).Value = 7000
$ws.Cells.Item(23, 12).Value = 7500
$ws.Cells.Item(23, 13).Value = 7250
$ws.Cells.Item(23, 15).Value = 'Provincia de Diguillín'
$ws.Cells.Item(23, 16).Value = 725
